$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.772.14'
$ws.Range('E2').Value = '  +3.78%  '
$ws.Range('D3').Value = '1.869.80'
$ws.Range('E3').Value = '  +3.13%  '
$ws.Range('E4').Value = '  +0.38%  '
$ws.Range('D5').Value = "'232.56"
$ws.Range('E5').Value = '  +3.05%  '
$ws.Range('E6').Value = '  +3.64%  '
$ws.Range('D8').Value = "'42.72"
$ws.Range('E8').Value = '  +11.32%  '
$ws.Range('D9').Value = "'0.312"
$ws.Range('E9').Value = '  +7.50%  '
$ws.Range('D10').Value = "'0.0699"
$ws.Range('E10').Value = '  +3.38%  '
$ws.Range('E11').Value = '  +4.23%  '
$ws.Range('D12').Value = '2.142.44'
$ws.Range('E12').Value = '  +3.24%  '
$ws.Range('D13').Value = "'11.74"
$ws.Range('E13').Value = '  +4.52%  '
$ws.Range('D14').Value = '1.875.65'
$ws.Range('E14').Value = '  +3.13%  '
$ws.Range('D15').Value = "'0.686"
$ws.Range('E15').Value = '  +8.31%  '
$ws.Range('D16').Value = "'4.78"
$ws.Range('E16').Value = '  +8.00%  '
$ws.Range('D17').Value = '35.797.83'
$ws.Range('E17').Value = '  +3.88%  '
$ws.Range('D18').Value = "'70.74"
$ws.Range('E18').Value = '  +3.53%  '
$ws.Range('E19').Value = '  +4.54%  '
$ws.Range('D20').Value = "'249.11"
$ws.Range('E20').Value = '  +2.42%  '
$ws.Range('D21').Value = "'12.51"
$ws.Range('E21').Value = '  +11.17%  '
$ws.Range('D22').Value = "'4.81"
$ws.Range('E22').Value = '  +16.57%  '
$ws.Range('E23').Value = '  +0.27%  '
$ws.Range('E24').Value = '  +1.84%  '
$ws.Range('D25').Value = "'171.82"
$ws.Range('E25').Value = '  +0.91%  '
$ws.Range('D26').Value = "'8.09"
$ws.Range('E26').Value = '  +3.51%  '
$ws.Range('D27').Value = "'18.00"
$ws.Range('E27').Value = '  +2.35%  '
$ws.Range('D28').Value = "'0.123"
$ws.Range('E28').Value = '  +2.03%  '
$ws.Range('E29').Value = '  +16.59%  '
$ws.Range('D30').Value = "'1.01"
$ws.Range('E30').Value = '  +0.42%  '
$ws.Range('D31').Value = '3.333.10'
$ws.Range('E31').Value = '  +37.18%  '
$ws.Range('D32').Value = "'0.0554"
$ws.Range('E32').Value = '  +6.85%  '
$ws.Range('D33').Value = "'3.98"
$ws.Range('E33').Value = '  +4.82%  '
$ws.Range('E34').Value = '  +6.56%  '
$ws.Range('E35').Value = '  +4.93%  '
$ws.Range('D36').Value = "'99.13"
$ws.Range('E36').Value = '  +21.17%  '
$ws.Range('D37').Value = "'0.692"
$ws.Range('E37').Value = '  +7.06%  '
$ws.Range('D38').Value = "'2.52"
$ws.Range('E38').Value = '  +6.88%  '
$ws.Range('D39').Value = '1.363.42'
$ws.Range('E39').Value = '  +0.21%  '
$ws.Range('E40').Value = '  +3.34%  '
$ws.Range('E41').Value = '  +5.96%  '
$ws.Range('E42').Value = '  +8.14%  '
$ws.Range('D43').Value = "'15.08"
$ws.Range('E43').Value = '  +9.35%  '
$ws.Range('E44').Value = '  +4.68%  '
$ws.Range('E45').Value = '  +1.62%  '
$ws.Range('E46').Value = '  +1.03%  '
$ws.Range('D47').Value = "'6.29"
$ws.Range('E47').Value = '  +9.07%  '
$ws.Range('E48').Value = '  +2.26%  '
$ws.Range('D49').Value = '2.040.30'
$ws.Range('E49').Value = '  +3.26%  '
$ws.Range('D50').Value = "'105.31"
$ws.Range('E50').Value = '  +3.18%  '
$ws.Range('E51').Value = '  +0.38%  '
